# Update "Name of Algo" result values in column B (imputed values)
# for the KNN algorithm result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.156000000000001
$ws.Range("B6").Value = 5.787
$ws.Range("B7").Value = 5.396000000000001
$ws.Range("B8").Value = 6.01
$ws.Range("B16").Value = 5.464
$ws.Range("B20").Value = 8.006
$ws.Range("B21").Value = 8.894
$ws.Range("B28").Value = 5.624000000000001
$ws.Range("B29").Value = 5.306
$ws.Range("B30").Value = 6.159000000000001
$ws.Range("B32").Value = 6.769999999999999
$ws.Range("B40").Value = 9.178999999999998
$ws.Range("B46").Value = 6.382000000000001
$ws.Range("B51").Value = 5.548
$ws.Range("B52").Value = 5.807
$ws.Range("B57").Value = 5.207
$ws.Range("B59").Value = 4.92
$ws.Range("B62").Value = 5.386
$ws.Range("B66").Value = 5.013
$ws.Range("B73").Value = 6.914000000000001
$ws.Range("B74").Value = 9.164999999999999
$ws.Range("B77").Value = 5.586
$ws.Range("B92").Value = 4.893000000000001
$ws.Range("B100").Value = 5.968
